# Applies "Changes for doc_is and chamber_id continues.." to the Query sheet.
# Adds a new "Investigation List" query row (row 9) plus four blank,
# formatted spacer rows (10-13) beneath it, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Query")

# --- Row 9: Investigation List query -------------------------------------
# Shared strings are allocated in the order the cells are written, so write
# column C (label) before column B (query) to reproduce the original index
# order: 47=Investigation List, 48=query($q%), 49=query example(FB%).

# Column C: label
$cC9 = $ws.Cells.Item(9, 3)
$cC9.Value = "Investigation List"
$cC9.Font.Name = "Trebuchet MS"
$cC9.Font.Size = 9

# Column B: the parameterised query
$cB9 = $ws.Cells.Item(9, 2)
$cB9.Value = "select a.investigation_name ,  a.ID  from  investigation_master a where a.investigation_name LIKE '`$q%' and STATUS = 'ACTIVE' AND a.chamber_id='`$chamber_name' AND a.doc_id='`$doc_name'"
$cB9.Font.Name = "Trebuchet MS"
$cB9.Font.Size = 9
$cB9.WrapText = $true

# Column D: worked example of the query
$cD9 = $ws.Cells.Item(9, 4)
$cD9.Value = "select a.investigation_name ,  a.ID  from  investigation_master a where a.investigation_name LIKE 'FB%' and STATUS = 'ACTIVE' AND a.chamber_id='sos' AND a.doc_id='sroy'"
$cD9.Font.Name = "Trebuchet MS"
$cD9.Font.Size = 9
$cD9.WrapText = $true

$ws.Rows.Item(9).RowHeight = 30

# --- Rows 10-13: blank spacer rows (formatted like column B wrap style) --
foreach ($r in 10..13) {
    $c = $ws.Cells.Item($r, 2)
    $c.Font.Name = "Trebuchet MS"
    $c.Font.Size = 9
    $c.WrapText = $true
}

# --- Update the selected / active cell ------------------------------------
$ws.Activate()
$ws.Range("D10").Select()
